$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation needs to be inserted as row 51 (pushing the
# existing rows 51-53 down to 52-54). Insert a blank row at position 51.
$ws.Rows.Item(51).Insert()

# Give the new row's date cell (D51) the same date number format used by the
# surrounding date cells (style index 2 / "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat

# Populate the newly inserted row with the new weekly data.
$ws.Cells.Item(51, 1).Value2 = 4
$ws.Cells.Item(51, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(51, 3).Value2 = "Los Lagos"
$ws.Cells.Item(51, 4).Value2 = 44509
$ws.Cells.Item(51, 5).Value2 = 10
$ws.Cells.Item(51, 6).Value2 = 100112026
$ws.Cells.Item(51, 7).Value2 = "Haba"
$ws.Cells.Item(51, 8).Value2 = "Sin especificar"
$ws.Cells.Item(51, 9).Value2 = "Primera"
$ws.Cells.Item(51, 10).Value2 = 160
$ws.Cells.Item(51, 11).Value2 = 10000
$ws.Cells.Item(51, 12).Value2 = 10000
$ws.Cells.Item(51, 13).Value2 = 10000
$ws.Cells.Item(51, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(51, 15).Value2 = "Región del Maule"
$ws.Cells.Item(51, 16).Value2 = 400
$ws.Cells.Item(51, 17).Value2 = 25
$ws.Cells.Item(51, 18).Value2 = "Hortaliza"
